{"js": "// Apply the Review_220 edit: update the daily-paper header/date, replace the\n// paper title + all four body paragraphs with the new \"grok\" review text, and\n// append the new link + sign-off paragraphs at the end of the body.\nconst newTexts = [\n  \"\u26a1\ufe0f\ud83d\ude80 \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 10.06.24:\u26a1\ufe0f\ud83d\ude80\",\n  \"Learning to grok: Emergence of in-context learning and skill composition in modular arithmetic tasks\",\n  \"\u05d0\u05d7\u05d3 \u05d4\u05ea\u05d5\u05e4\u05e2\u05d5\u05ea \u05d4\u05de\u05e8\u05ea\u05e7\u05d5\u05ea \u05d1\u05dc\u05de\u05d9\u05d3\u05d4 \u05e2\u05de\u05d5\u05e7\u05d4 \u05d4\u05d9\u05d0 \u05d2\u05e8\u05d5\u05e7\u05d9\u05e0\u05d2 - \u05e9\u05d4\u05d9\u05d0 \u05de\u05e2\u05d1\u05e8 \u05f4\u05e4\u05ea\u05d0\u05d5\u05de\u05d9\u05f4 \u05e9\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e2\u05de\u05d5\u05e7\u05d5\u05ea \u05dc\u05de\u05e6\u05d1 \u05e9\u05dc \u05d4\u05db\u05dc\u05dc\u05d4 \u05de\u05d4\u05de\u05e6\u05d1 \u05e9\u05dc overfitting \u05dc\u05de\u05e9\u05dc \u05d0\u05d7\u05e8\u05d9 \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d0\u05d5\u05d3 \u05d0\u05e8\u05d5\u05da. \u05d4\u05e8\u05d9 \u05d9\u05d3\u05d5\u05e2 \u05e9\u05d0\u05dd \u05e2\u05d1\u05d5\u05e8 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e0\u05ea\u05d5\u05df \u05d5\u05e8\u05e9\u05ea \u05e2\u05de\u05d5\u05e7\u05d4 \u05d1\u05e2\u05dc\u05ea \u05d9\u05db\u05d5\u05dc\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d2\u05d1\u05d5\u05d4\u05d4 \u05de\u05e1\u05e4\u05d9\u05e7 (representativeness) \u05d0\u05d7\u05e8\u05d9 \u05e9\u05dc\u05d1 \u05de\u05e1\u05d5\u05d9\u05dd \u05d1\u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05e0\u05d5 \u05e0\u05d2\u05d9\u05e2 \u05dc-overfitting \u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05de\u05e6\u05d1 \u05e9\u05d1\u05d5 \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05d9\u05d9\u05dc\u05db\u05d5 \u05d5\u05d9\u05e9\u05ea\u05e4\u05e8\u05d5 \u05e2\u05d1\u05d5\u05e8 \u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05d5\u05dc\u05dd \u05d4\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05e2\u05dc \u05e1\u05d8 \u05d4\u05d5\u05dc\u05d9\u05d3\u05e6\u05d9\u05d4 \u05d9\u05e1\u05e4\u05d2\u05d5 \u05d9\u05e8\u05d9\u05d3\u05d4 \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd. \",\n  \"\u05de\u05d4 \u05e9\u05de\u05d2\u05e0\u05d9\u05d1 \u05d5\u05de\u05e4\u05ea\u05d9\u05e2 \u05d1\u05d2\u05e8\u05d5\u05e7\u05d9\u05e0\u05d2 \u05e9\u05e2\u05d1\u05d5\u05e8 \u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05e8\u05d5\u05da \u05de\u05e1\u05e4\u05d9\u05e7 \u05de\u05d2\u05d9\u05e2 \u05d4\u05de\u05e6\u05d1 \u05e9\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05dc \u05e1\u05d8 \u05d4\u05d5\u05d5\u05dc\u05d9\u05d3\u05e6\u05d9\u05d4 \u05de\u05ea\u05d7\u05d9\u05dc\u05d9\u05dd \u05dc\u05e2\u05dc\u05d5\u05ea \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d0\u05dc\u05d4 \u05e2\u05dc \u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05d2\u05d9\u05e2 \u05dc\u05e9\u05dc\u05d1 \u05e9\u05dc \u05d4\u05db\u05dc\u05dc\u05d4 \u05d0\u05de\u05d9\u05ea\u05d9\u05ea. \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e9\u05ea\u05d5\u05e4\u05e2\u05d4 \u05d3\u05d5\u05de\u05d4 \u05de\u05ea\u05e8\u05d7\u05e9\u05ea \u05d1\u05ea\u05e0\u05d0\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05de\u05d9\u05dd \u05d0\u05dd \u05d0\u05e0\u05d5 \u05de\u05d2\u05d3\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05e7\u05d9\u05d1\u05d5\u05dc\u05ea \u05d4\u05de\u05d5\u05d3\u05dc (\u05de\u05e1' \u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd) \u05db\u05d0\u05e9\u05e8 \u05d2\u05d5\u05d3\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d5\u05de\u05e9\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e0\u05d5\u05ea\u05e8\u05d9\u05dd \u05e7\u05d1\u05d5\u05e2\u05d9\u05dd) \u05d5\u05d2\u05dd \u05db\u05d0\u05e9\u05e8 \u05d0\u05e0\u05d5 \u05de\u05d2\u05d3\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05d2\u05d5\u05d3\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05ea\u05d5\u05da \u05e9\u05de\u05d9\u05e8\u05d4 \u05e9\u05dc \u05de\u05e9\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e7\u05d9\u05d1\u05d5\u05dc\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05e7\u05d1\u05d5\u05e2\u05d9\u05dd. \",\n  \"\u05dc\u05de\u05e2\u05e9\u05d4 \u05ea\u05d5\u05e4\u05e2\u05d5\u05ea \u05d0\u05dc\u05d5 \u05e9\u05d9\u05d9\u05db\u05d5\u05ea \u05dc\u05de\u05e9\u05e4\u05d7\u05ea double descent (\u05d9\u05e9 \u05d2\u05dd multiple descent) \u05e9\u05e0\u05d7\u05e7\u05e8\u05d4 \u05e8\u05d1\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d7\u05d5\u05e7\u05e8 \u05d3\u05d2\u05d5\u05dc \u05de\u05d9\u05e9\u05d4 \u05d1\u05dc\u05e7\u05d9\u05df. \u05d4\u05ea\u05d5\u05e4\u05e2\u05d4 \u05e2\u05e6\u05de\u05d4 \u05e0\u05ea\u05d2\u05dc\u05ea\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05d9\u05d5\u05ea\u05e8 \u05de- 30 \u05e9\u05e0\u05d4 (\u05de\u05d9 \u05e9\u05e8\u05d5\u05e6\u05d4 \u05dc\u05d4\u05ea\u05e2\u05de\u05e7 \u05d1\u05e0\u05d5\u05e9\u05d0 \u05ea\u05e2\u05e7\u05d1\u05d5 \u05d0\u05d7\u05e8\u05d9 https://www.linkedin.com/in/charlesmartin14/ - \u05d4\u05d5\u05d0 \u05d0\u05d7\u05d3 \u05d4\u05de\u05d5\u05de\u05d7\u05d9\u05dd \u05d4\u05d2\u05d3\u05d5\u05dc\u05d9\u05dd).\",\n  \"\u05d0\u05d5\u05e7\u05d9\u05d9, \u05d0\u05d6 \u05de\u05d4 \u05e2\u05e9\u05d4 \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05e0\u05e1\u05e7\u05e8? \u05d4\u05d5\u05d0 \u05d7\u05e7\u05e8 \u05ea\u05d5\u05e4\u05e2\u05ea \u05d2\u05e8\u05d5\u05e7\u05d9\u05e0\u05d2 \u05db\u05d0\u05e9\u05e8 \u05de\u05ea\u05e8\u05d7\u05e9\u05ea \u05d0\u05dd \u05de\u05d2\u05d3\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05de\u05e1\u05e4\u05e8 \u05d4\u05de\u05e9\u05d9\u05de\u05d5\u05ea (\u05db\u05dc \u05de\u05e9\u05d9\u05de\u05d4 \u05d4\u05d9\u05d0 \u05e1\u05d5\u05d2 \u05e9\u05dc \u05e8\u05d2\u05e8\u05e1\u05d9\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d1\u05e9\u05d3\u05d4 \u05d4\u05de\u05d5\u05d3\u05d5\u05dc\u05d5(\u05e9\u05d0\u05e8\u05d9\u05ea)) \u05e9\u05e2\u05d1\u05d5\u05e8\u05df \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc (\u05db\u05de\u05d5\u05d1\u05df \u05dc\u05e7\u05d7\u05d5 \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4). \u05de\u05ea\u05d1\u05e8\u05e8 \u05db\u05d9 \u05d9\u05e9 \u05db\u05de\u05d4 \u05de\u05e9\u05d8\u05e8\u05d9\u05dd (\u05de\u05d5\u05d3\u05d9\u05dd) \u05e9\u05dc \u05d9\u05db\u05d5\u05dc\u05ea \u05d4\u05db\u05dc\u05dc\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05db\u05d0\u05e9\u05e8 \u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05e2\u05dd \u05d4\u05d9\u05d7\u05e1 \u05e9\u05dc \u05de\u05e1\u05e4\u05e8 \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e4\u05e8 \u05de\u05e9\u05d9\u05de\u05d4 \u05d5\u05e2\u05dd \u05de\u05e1\u05e4\u05e8 \u05d4\u05de\u05e9\u05d9\u05de\u05d4. \u05d1\u05d2\u05d3\u05d5\u05dc \u05de\u05d0\u05d5\u05d3 \u05d0\u05dd \u05e0\u05d5\u05ea\u05e0\u05d9\u05dd \u05de\u05e1\u05e4\u05d9\u05e7 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d2\u05d3\u05d5\u05dc \u05de\u05e1\u05e4\u05d9\u05e7 \u05d5\u05de\u05e1\u05e4\u05e8 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e4\u05e8 \u05de\u05e9\u05d9\u05de\u05d4 \u05d2\u05d3\u05d5\u05dc \u05de\u05e1\u05e4\u05d9\u05e7 \u05d0\u05d6 \u05de\u05d2\u05d9\u05e2\u05d9\u05dd \u05dc\u05d4\u05db\u05dc\u05dc\u05d4 \u05d0\u05de\u05d9\u05ea\u05d9\u05ea \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05d0\u05db\u05df \u05dc\u05d5\u05de\u05d3 \u05d0\u05ea \u05d4\u05de\u05e9\u05d9\u05de\u05d4 \u05d1\u05de\u05dc\u05d5\u05d0\u05d4).\",\n  \"https://arxiv.org/abs/2406.02550\",\n  \"\u05e7\u05e8\u05d9\u05d0\u05d4 \u05de\u05d4\u05e0\u05d4!\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The original document has exactly 6 paragraphs; the first 6 entries of\n// newTexts replace them 1:1, and the trailing 2 entries are brand-new\n// paragraphs appended after the (now-updated) final paragraph.\nconst existingCount = paragraphs.items.length;\n\nfor (let i = 0; i < existingCount && i < newTexts.length; i++) {\n  paragraphs.items[i].insertText(newTexts[i], \"Replace\");\n}\nawait context.sync();\n\nlet insertAfter = paragraphs.items[existingCount - 1];\nfor (let i = existingCount; i < newTexts.length; i++) {\n  insertAfter = insertAfter.insertParagraph(newTexts[i], \"After\");\n}\nawait context.sync();\n", "ps1": "# Apply the Review_220 edit via the Word COM object model: update the\n# daily-paper header/date, replace the paper title and the four body\n# paragraphs with the new \"grok\" review text, then append the new link\n# and sign-off paragraphs at the end of the document.\n$d = $word.ActiveDocument\n\n$newTexts = @(\n  '\u26a1\ufe0f\ud83d\ude80 \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 10.06.24:\u26a1\ufe0f\ud83d\ude80',\n  'Learning to grok: Emergence of in-context learning and skill composition in modular arithmetic tasks',\n  '\u05d0\u05d7\u05d3 \u05d4\u05ea\u05d5\u05e4\u05e2\u05d5\u05ea \u05d4\u05de\u05e8\u05ea\u05e7\u05d5\u05ea \u05d1\u05dc\u05de\u05d9\u05d3\u05d4 \u05e2\u05de\u05d5\u05e7\u05d4 \u05d4\u05d9\u05d0 \u05d2\u05e8\u05d5\u05e7\u05d9\u05e0\u05d2 - \u05e9\u05d4\u05d9\u05d0 \u05de\u05e2\u05d1\u05e8 \u05f4\u05e4\u05ea\u05d0\u05d5\u05de\u05d9\u05f4 \u05e9\u05dc \u05e8\u05e9\u05ea\u05d5\u05ea \u05e2\u05de\u05d5\u05e7\u05d5\u05ea \u05dc\u05de\u05e6\u05d1 \u05e9\u05dc \u05d4\u05db\u05dc\u05dc\u05d4 \u05de\u05d4\u05de\u05e6\u05d1 \u05e9\u05dc overfitting \u05dc\u05de\u05e9\u05dc \u05d0\u05d7\u05e8\u05d9 \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d0\u05d5\u05d3 \u05d0\u05e8\u05d5\u05da. \u05d4\u05e8\u05d9 \u05d9\u05d3\u05d5\u05e2 \u05e9\u05d0\u05dd \u05e2\u05d1\u05d5\u05e8 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e0\u05ea\u05d5\u05df \u05d5\u05e8\u05e9\u05ea \u05e2\u05de\u05d5\u05e7\u05d4 \u05d1\u05e2\u05dc\u05ea \u05d9\u05db\u05d5\u05dc\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d2\u05d1\u05d5\u05d4\u05d4 \u05de\u05e1\u05e4\u05d9\u05e7 (representativeness) \u05d0\u05d7\u05e8\u05d9 \u05e9\u05dc\u05d1 \u05de\u05e1\u05d5\u05d9\u05dd \u05d1\u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05e0\u05d5 \u05e0\u05d2\u05d9\u05e2 \u05dc-overfitting \u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05de\u05e6\u05d1 \u05e9\u05d1\u05d5 \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05d9\u05d9\u05dc\u05db\u05d5 \u05d5\u05d9\u05e9\u05ea\u05e4\u05e8\u05d5 \u05e2\u05d1\u05d5\u05e8 \u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05d5\u05dc\u05dd \u05d4\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05e2\u05dc \u05e1\u05d8 \u05d4\u05d5\u05dc\u05d9\u05d3\u05e6\u05d9\u05d4 \u05d9\u05e1\u05e4\u05d2\u05d5 \u05d9\u05e8\u05d9\u05d3\u05d4 \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd. ',\n  '\u05de\u05d4 \u05e9\u05de\u05d2\u05e0\u05d9\u05d1 \u05d5\u05de\u05e4\u05ea\u05d9\u05e2 \u05d1\u05d2\u05e8\u05d5\u05e7\u05d9\u05e0\u05d2 \u05e9\u05e2\u05d1\u05d5\u05e8 \u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05e8\u05d5\u05da \u05de\u05e1\u05e4\u05d9\u05e7 \u05de\u05d2\u05d9\u05e2 \u05d4\u05de\u05e6\u05d1 \u05e9\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05dc \u05e1\u05d8 \u05d4\u05d5\u05d5\u05dc\u05d9\u05d3\u05e6\u05d9\u05d4 \u05de\u05ea\u05d7\u05d9\u05dc\u05d9\u05dd \u05dc\u05e2\u05dc\u05d5\u05ea \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d0\u05dc\u05d4 \u05e2\u05dc \u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05d2\u05d9\u05e2 \u05dc\u05e9\u05dc\u05d1 \u05e9\u05dc \u05d4\u05db\u05dc\u05dc\u05d4 \u05d0\u05de\u05d9\u05ea\u05d9\u05ea. \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e9\u05ea\u05d5\u05e4\u05e2\u05d4 \u05d3\u05d5\u05de\u05d4 \u05de\u05ea\u05e8\u05d7\u05e9\u05ea \u05d1\u05ea\u05e0\u05d0\u05d9\u05dd \u05de\u05e1\u05d5\u05d9\u05de\u05d9\u05dd \u05d0\u05dd \u05d0\u05e0\u05d5 \u05de\u05d2\u05d3\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05e7\u05d9\u05d1\u05d5\u05dc\u05ea \u05d4\u05de\u05d5\u05d3\u05dc (\u05de\u05e1'' \u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd) \u05db\u05d0\u05e9\u05e8 \u05d2\u05d5\u05d3\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d5\u05de\u05e9\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e0\u05d5\u05ea\u05e8\u05d9\u05dd \u05e7\u05d1\u05d5\u05e2\u05d9\u05dd) \u05d5\u05d2\u05dd \u05db\u05d0\u05e9\u05e8 \u05d0\u05e0\u05d5 \u05de\u05d2\u05d3\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05d2\u05d5\u05d3\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05ea\u05d5\u05da \u05e9\u05de\u05d9\u05e8\u05d4 \u05e9\u05dc \u05de\u05e9\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e7\u05d9\u05d1\u05d5\u05dc\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05e7\u05d1\u05d5\u05e2\u05d9\u05dd. ',\n  '\u05dc\u05de\u05e2\u05e9\u05d4 \u05ea\u05d5\u05e4\u05e2\u05d5\u05ea \u05d0\u05dc\u05d5 \u05e9\u05d9\u05d9\u05db\u05d5\u05ea \u05dc\u05de\u05e9\u05e4\u05d7\u05ea double descent (\u05d9\u05e9 \u05d2\u05dd multiple descent) \u05e9\u05e0\u05d7\u05e7\u05e8\u05d4 \u05e8\u05d1\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d7\u05d5\u05e7\u05e8 \u05d3\u05d2\u05d5\u05dc \u05de\u05d9\u05e9\u05d4 \u05d1\u05dc\u05e7\u05d9\u05df. \u05d4\u05ea\u05d5\u05e4\u05e2\u05d4 \u05e2\u05e6\u05de\u05d4 \u05e0\u05ea\u05d2\u05dc\u05ea\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05d9\u05d5\u05ea\u05e8 \u05de- 30 \u05e9\u05e0\u05d4 (\u05de\u05d9 \u05e9\u05e8\u05d5\u05e6\u05d4 \u05dc\u05d4\u05ea\u05e2\u05de\u05e7 \u05d1\u05e0\u05d5\u05e9\u05d0 \u05ea\u05e2\u05e7\u05d1\u05d5 \u05d0\u05d7\u05e8\u05d9 https://www.linkedin.com/in/charlesmartin14/ - \u05d4\u05d5\u05d0 \u05d0\u05d7\u05d3 \u05d4\u05de\u05d5\u05de\u05d7\u05d9\u05dd \u05d4\u05d2\u05d3\u05d5\u05dc\u05d9\u05dd).',\n  '\u05d0\u05d5\u05e7\u05d9\u05d9, \u05d0\u05d6 \u05de\u05d4 \u05e2\u05e9\u05d4 \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05e0\u05e1\u05e7\u05e8? \u05d4\u05d5\u05d0 \u05d7\u05e7\u05e8 \u05ea\u05d5\u05e4\u05e2\u05ea \u05d2\u05e8\u05d5\u05e7\u05d9\u05e0\u05d2 \u05db\u05d0\u05e9\u05e8 \u05de\u05ea\u05e8\u05d7\u05e9\u05ea \u05d0\u05dd \u05de\u05d2\u05d3\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05de\u05e1\u05e4\u05e8 \u05d4\u05de\u05e9\u05d9\u05de\u05d5\u05ea (\u05db\u05dc \u05de\u05e9\u05d9\u05de\u05d4 \u05d4\u05d9\u05d0 \u05e1\u05d5\u05d2 \u05e9\u05dc \u05e8\u05d2\u05e8\u05e1\u05d9\u05d4 \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea \u05d1\u05e9\u05d3\u05d4 \u05d4\u05de\u05d5\u05d3\u05d5\u05dc\u05d5(\u05e9\u05d0\u05e8\u05d9\u05ea)) \u05e9\u05e2\u05d1\u05d5\u05e8\u05df \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc (\u05db\u05de\u05d5\u05d1\u05df \u05dc\u05e7\u05d7\u05d5 \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4). \u05de\u05ea\u05d1\u05e8\u05e8 \u05db\u05d9 \u05d9\u05e9 \u05db\u05de\u05d4 \u05de\u05e9\u05d8\u05e8\u05d9\u05dd (\u05de\u05d5\u05d3\u05d9\u05dd) \u05e9\u05dc \u05d9\u05db\u05d5\u05dc\u05ea \u05d4\u05db\u05dc\u05dc\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05db\u05d0\u05e9\u05e8 \u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05e2\u05dd \u05d4\u05d9\u05d7\u05e1 \u05e9\u05dc \u05de\u05e1\u05e4\u05e8 \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e4\u05e8 \u05de\u05e9\u05d9\u05de\u05d4 \u05d5\u05e2\u05dd \u05de\u05e1\u05e4\u05e8 \u05d4\u05de\u05e9\u05d9\u05de\u05d4. \u05d1\u05d2\u05d3\u05d5\u05dc \u05de\u05d0\u05d5\u05d3 \u05d0\u05dd \u05e0\u05d5\u05ea\u05e0\u05d9\u05dd \u05de\u05e1\u05e4\u05d9\u05e7 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d2\u05d3\u05d5\u05dc \u05de\u05e1\u05e4\u05d9\u05e7 \u05d5\u05de\u05e1\u05e4\u05e8 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e4\u05e8 \u05de\u05e9\u05d9\u05de\u05d4 \u05d2\u05d3\u05d5\u05dc \u05de\u05e1\u05e4\u05d9\u05e7 \u05d0\u05d6 \u05de\u05d2\u05d9\u05e2\u05d9\u05dd \u05dc\u05d4\u05db\u05dc\u05dc\u05d4 \u05d0\u05de\u05d9\u05ea\u05d9\u05ea \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05d0\u05db\u05df \u05dc\u05d5\u05de\u05d3 \u05d0\u05ea \u05d4\u05de\u05e9\u05d9\u05de\u05d4 \u05d1\u05de\u05dc\u05d5\u05d0\u05d4).',\n  'https://arxiv.org/abs/2406.02550',\n  '\u05e7\u05e8\u05d9\u05d0\u05d4 \u05de\u05d4\u05e0\u05d4!'\n)\n\n$existingCount = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $existingCount -and $i -le $newTexts.Count; $i++) {\n  $d.Paragraphs($i).Range.Text = $newTexts[$i - 1]\n}\n\n$tail = $d.Paragraphs($existingCount).Range\nfor ($i = $existingCount + 1; $i -le $newTexts.Count; $i++) {\n  $tail.InsertParagraphAfter()\n  $tail = $d.Paragraphs($i).Range\n  $tail.Text = $newTexts[$i - 1]\n}\n\n"}
